# Master Data Tables - Test Data / master-user_detail.xlsx
# "Updated Master data as per 16th May Refresh"
# Appends 3 new user rows (Nikola Tesla, Graham Bell, Albert Miles) to Sheet1,
# following the same shape as the existing rows (row 33 is the last existing
# data row, so the new rows land at 34-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New users to append, in order.
$newUsers = @(
    @{ Id = 110033; Uin = 9317596771; Name = "Nikola Tesla"; Email = "nikola.tesla@xyz.com";  Mobile = 818876434 },
    @{ Id = 110034; Uin = 9317596772; Name = "Graham Bell";  Email = "graham.bell@xyz.com";   Mobile = 818876435 },
    @{ Id = 110035; Uin = 9317596773; Name = "Albert Miles"; Email = "albert.miles@xyz.com";  Mobile = 818876436 }
)

$startRow = 34

# Fill column-by-column (id, uin, name, email, mobile, ...) across all new
# rows so any newly-introduced shared strings land in the same order as the
# source data (all names first, then all emails), matching the name/email
# columns being populated as a block.
$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 1).Value = $u.Id
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 2).Value = $u.Uin
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 3).Value = $u.Name
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 4).Value = $u.Email
    # Column D (email) uses the same applied style as the rest of the email
    # column in the existing data.
    $ws.Cells.Item($row, 4).HorizontalAlignment = -4131
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 5).Value = $u.Mobile
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 6).Value = "ACT"
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 7).Value = "eng"
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 8).Value = "PWD"
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 9).Value = $true
    # Column I (is_active) matches the existing left-aligned boolean style.
    $ws.Cells.Item($row, 9).HorizontalAlignment = -4131
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 10).Value = "superadmin"
    $row = $row + 1
}

$row = $startRow
foreach ($u in $newUsers) {
    $ws.Cells.Item($row, 11).Value = "now()"
    $row = $row + 1
}

# Move the selection below the newly added rows, mirroring the author's
# cursor position after entering the new data (select from the row after the
# last data row down to the end of the sheet).
$lastRow = $row
$ws.Range("A" + $lastRow + ":XFD1048576").Select()
